# Automatische test-sync: 2025-07-27 19:25:50
#
# Appends a new log entry (Testmail #7) to the "Logs" sheet, extends the
# conditional-formatting ranges to cover the new row, and updates the
# "Dashboard" summary sheet (category counts / order) to reflect the new
# "Productinformatie" entry.

$wb = $excel.ActiveWorkbook

# --- 1. Append the new row to the Logs sheet -------------------------------
$logs = $wb.Worksheets.Item("Logs")

$newRow = 9
$logs.Cells.Item($newRow, 1).Value = "Is dit artikel nog op voorraad?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #7: Is dit artikel nog op voorraad?"
$logs.Cells.Item($newRow, 4).Value = "Productinformatie"
$logs.Cells.Item($newRow, 5).Value = "FALLBACK_BLOCKED"
$logs.Cells.Item($newRow, 6).Value = "2025-07-27 19:25:33"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Nee"
$logs.Cells.Item($newRow, 9).Value = "Ja"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# --- 2. Extend the conditional formatting ranges to include the new row ---
$logs.Range("D2:D8").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D9"))
$logs.Range("G2:G8").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G9"))
$logs.Range("H2:H8").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H9"))
$logs.Range("I2:I8").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I9"))
$logs.Range("J2:J8").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J9"))

# --- 3. Update the Dashboard summary sheet ---------------------------------
# "Productinformatie" now has 2 occurrences and moves above
# "Bestelling / Levering" (which stays at 1).
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Cells.Item(4, 1).Value = "Productinformatie"
$dashboard.Cells.Item(4, 2).Value = 2

$dashboard.Cells.Item(5, 1).Value = "Bestelling / Levering"
$dashboard.Cells.Item(5, 2).Value = 1
